# Generate Report for Handback
# Updates the localization-status workbook to reflect that the de-de and
# zh-cn handback packages have been generated: status moves from
# "Ready for handoff" to "Handed back: in sync with en-US", the target/
# handback file names + handback datetime are filled in, and the
# corresponding columns are widened to fit the new (longer) content.

$wb = $excel.ActiveWorkbook

$mdFileName = "8334310f-08ad-4762-bbb0-ead26ddd535c.md"
$mdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/476d99cbf29b0af1d3716c8f97a70ae3aefb5bf4/e2e/8334310f-08ad-4762-bbb0-ead26ddd535c.md"
$statusText = "Handed back: in sync with en-US"

# Hyperlink font formatting (matches the existing "HyperLink" cell style:
# underlined, Calibri 11, color #6495ED) used for newly hyperlinked cells.
$hyperlinkColor = 15570276

# ---------------------------------------------------------------------
# Overview sheet: update Status cells (E2, F2) and widen their columns
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Columns.Item(5).ColumnWidth = 29.2
$wsOverview.Columns.Item(6).ColumnWidth = 29.2

# ---------------------------------------------------------------------
# zh-cn sheet: status, target/handback file + datetime, widen columns
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $statusText
$wsZh.Columns.Item(3).ColumnWidth = 29.2

$wsZh.Range("I2").Value = $mdFileName
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $mdUrl, [Type]::Missing, [Type]::Missing, $mdFileName)
$wsZh.Range("I2").Style = "HyperLink"
$wsZh.Range("I2").Font.Underline = 2
$wsZh.Range("I2").Font.Color = $hyperlinkColor

$wsZh.Range("J2").Value = "8334310f-08ad-4762-bbb0-ead26ddd535c.8f595f0fc4d32ff2d9a7ba2aad79a008fcbb94b4.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-09-04 23:04:09"

$wsZh.Columns.Item(9).ColumnWidth = 39.2
$wsZh.Columns.Item(10).ColumnWidth = 39.2

# ---------------------------------------------------------------------
# de-de sheet: status, target/handback file + datetime, widen columns
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $statusText
$wsDe.Columns.Item(3).ColumnWidth = 29.2

$wsDe.Range("I2").Value = $mdFileName
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $mdUrl, [Type]::Missing, [Type]::Missing, $mdFileName)
$wsDe.Range("I2").Style = "HyperLink"
$wsDe.Range("I2").Font.Underline = 2
$wsDe.Range("I2").Font.Color = $hyperlinkColor

$wsDe.Range("J2").Value = "8334310f-08ad-4762-bbb0-ead26ddd535c.8f595f0fc4d32ff2d9a7ba2aad79a008fcbb94b4.de-de.xlf"
$wsDe.Range("K2").Value = "2016-09-04 23:04:16"

$wsDe.Columns.Item(9).ColumnWidth = 39.2
$wsDe.Columns.Item(10).ColumnWidth = 39.2
